$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a brand-new row at row 12 (the most recent entry, 2022-01-27) ---
# This pushes the existing rows 12-19 down to rows 13-20.
$ws.Rows("12:12").Insert()

$ws.Range("A12").Value = 3
$ws.Range("B12").Value = "Femacal de La Calera"
$ws.Range("C12").Value = "Coquimbo"
$ws.Range("D12").Value = 44588
$ws.Range("D12").NumberFormat = $ws.Range("D13").NumberFormat
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100104
$ws.Range("H12").Value = "Frutos de pepita"
$ws.Range("I12").Value = 100104001
$ws.Range("J12").Value = "Granada"
$ws.Range("K12").Value = "Wonderfull"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 85
$ws.Range("N12").Value = 19000
$ws.Range("O12").Value = 20000
$ws.Range("P12").Value = 19529
$ws.Range("Q12").Value = "$/caja 14 kilos granel"
$ws.Range("R12").Value = "Provincia de Limarí"
$ws.Range("S12").Value = 1395
$ws.Range("T12").Value = 14

# --- Insert a second brand-new row at row 20 (2022-01-24) ---
# At this point rows 13-20 hold the old rows 12-19, so the old last row
# (originally row 19, now row 20) is pushed down to row 21.
$ws.Rows("20:20").Insert()

$ws.Range("A20").Value = 3
$ws.Range("B20").Value = "Femacal de La Calera"
$ws.Range("C20").Value = "Coquimbo"
$ws.Range("D20").Value = 44585
$ws.Range("D20").NumberFormat = $ws.Range("D19").NumberFormat
$ws.Range("E20").Value = 5
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100104
$ws.Range("H20").Value = "Frutos de pepita"
$ws.Range("I20").Value = 100104001
$ws.Range("J20").Value = "Granada"
$ws.Range("K20").Value = "Wonderfull"
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 50
$ws.Range("N20").Value = 22500
$ws.Range("O20").Value = 22500
$ws.Range("P20").Value = 22500
$ws.Range("Q20").Value = "$/caja 15 kilos empedrada"
$ws.Range("R20").Value = "Provincia de Limarí"
$ws.Range("S20").Value = 1500
$ws.Range("T20").Value = 15
